$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1.44
$ws.Range("E3").Value = 1.26
$ws.Range("C4").Value = 1.41
$ws.Range("C5").Value = 1.37
$ws.Range("D6").Value = 1.53
$ws.Range("G6").Value = 1.02
